# Update leve-profit calculation inputs/outputs per latest market data (scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 87.5
$ws.Range("I9").Value = 150
$ws.Range("K9").Value = 150
$ws.Range("M9").Value = 19

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2536.3635
$ws.Range("J29").Value = 2760
$ws.Range("L29").Value = 8280
$ws.Range("N29").Value = -8842

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3284.7144
$ws.Range("I74").Value = 2998.6
$ws.Range("K74").Value = 2998.6
$ws.Range("M74").Value = -2062.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3284.7144
$ws.Range("I77").Value = 2998.6
$ws.Range("K77").Value = 14993
$ws.Range("M77").Value = -10313

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5377.778
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 7316.6665
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 7316.6665
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -9562.666499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5377.778
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 7316.6665
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 36583.3325
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -47815.3325

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 649.85
$ws.Range("I92").Value = 527.6111
$ws.Range("K92").Value = 527.6111
$ws.Range("M92").Value = 720.3889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 30624
$ws.Range("J95").Value = 30624
$ws.Range("L95").Value = 30624
$ws.Range("N95").Value = -36116

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3815.65
$ws.Range("I116").Value = 3428.6428
$ws.Range("J116").Value = 4718.6665
$ws.Range("K116").Value = 3428.6428
$ws.Range("L116").Value = 4718.6665
$ws.Range("M116").Value = 13.35719999999992
$ws.Range("N116").Value = -11602.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 18527362
$ws.Range("I132").Value = 37051788
$ws.Range("K132").Value = 111155364
$ws.Range("M132").Value = -111152834

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 45455456
$ws.Range("I135").Value = 440.73334
$ws.Range("K135").Value = 3966.60006
$ws.Range("M135").Value = -1431.60006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 27174.285
$ws.Range("J136").Value = 29203.334
$ws.Range("L136").Value = 29203.334
$ws.Range("N136").Value = -39403.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3378.5405
$ws.Range("I32").Value = 3090.0317
$ws.Range("K32").Value = 3090.0317
$ws.Range("M32").Value = -2803.0317

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 10000
$ws.Range("L56").Value = 10000
$ws.Range("N56").Value = -11484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1597.4375
$ws.Range("I74").Value = 957.36365
$ws.Range("K74").Value = 957.36365
$ws.Range("M74").Value = -83.36365000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1597.4375
$ws.Range("I77").Value = 957.36365
$ws.Range("K77").Value = 4786.81825
$ws.Range("M77").Value = -418.8182500000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 149
$ws.Range("I11").Value = 119
$ws.Range("J11").Value = 179
$ws.Range("K11").Value = 119
$ws.Range("L11").Value = 179
$ws.Range("M11").Value = 21
$ws.Range("N11").Value = -459

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1754.4546
$ws.Range("I20").Value = 1199.8572
$ws.Range("J20").Value = 2725
$ws.Range("K20").Value = 1199.8572
$ws.Range("L20").Value = 2725
$ws.Range("M20").Value = -952.8571999999999
$ws.Range("N20").Value = -3219

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1000.7273
$ws.Range("I94").Value = 787.2
$ws.Range("K94").Value = 787.2
$ws.Range("M94").Value = -336.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4037.7
$ws.Range("I132").Value = 3875.8
$ws.Range("K132").Value = 11627.4
$ws.Range("M132").Value = -9097.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 358.8
$ws.Range("I13").Value = 199.75
$ws.Range("J13").Value = 995
$ws.Range("K13").Value = 599.25
$ws.Range("L13").Value = 2985
$ws.Range("M13").Value = -431.25
$ws.Range("N13").Value = -3321

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10206061
$ws.Range("I131").Value = 166667090
$ws.Range("J131").Value = 2081.2283
$ws.Range("K131").Value = 500001270
$ws.Range("L131").Value = 6243.6849
$ws.Range("M131").Value = -499996230
$ws.Range("N131").Value = -16323.6849

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 75002296
$ws.Range("I70").Value = 83335336
$ws.Range("K70").Value = 83335336
$ws.Range("M70").Value = -83335066

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 75002296
$ws.Range("I73").Value = 83335336
$ws.Range("K73").Value = 83335336
$ws.Range("M73").Value = -83334400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2917.647
$ws.Range("I80").Value = 2560
$ws.Range("J80").Value = 3168
$ws.Range("K80").Value = 2560
$ws.Range("L80").Value = 3168
$ws.Range("M80").Value = -1562
$ws.Range("N80").Value = -5164

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2917.647
$ws.Range("I83").Value = 2560
$ws.Range("J83").Value = 3168
$ws.Range("K83").Value = 12800
$ws.Range("L83").Value = 15840
$ws.Range("M83").Value = -7808
$ws.Range("N83").Value = -25824

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 866.6667
$ws.Range("I97").Value = 950
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 950
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -454
$ws.Range("N97").Value = -1692

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1865
$ws.Range("I122").Value = 1898.75
$ws.Range("J122").Value = 1730
$ws.Range("K122").Value = 5696.25
$ws.Range("L122").Value = 5190
$ws.Range("M122").Value = -3246.25
$ws.Range("N122").Value = -10090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1830.125
$ws.Range("J22").Value = 1830.125
$ws.Range("L22").Value = 1830.125
$ws.Range("N22").Value = -2420.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1830.125
$ws.Range("J27").Value = 1830.125
$ws.Range("L27").Value = 1830.125
$ws.Range("N27").Value = -2044.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1187.35
$ws.Range("I61").Value = 855.9231
$ws.Range("J61").Value = 1802.8572
$ws.Range("K61").Value = 855.9231
$ws.Range("L61").Value = 1802.8572
$ws.Range("M61").Value = -653.9231
$ws.Range("N61").Value = -2206.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1187.35
$ws.Range("I113").Value = 855.9231
$ws.Range("J113").Value = 1802.8572
$ws.Range("K113").Value = 855.9231
$ws.Range("L113").Value = 1802.8572
$ws.Range("M113").Value = 1314.0769
$ws.Range("N113").Value = -6142.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 23777.445
$ws.Range("I132").Value = 1024.3572
$ws.Range("K132").Value = 3073.0716
$ws.Range("M132").Value = -543.0715999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 48000
$ws.Range("J101").Value = 48000
$ws.Range("L101").Value = 48000
$ws.Range("N101").Value = -54490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 11900
$ws.Range("J103").Value = 11900
$ws.Range("L103").Value = 11900
$ws.Range("N103").Value = -14244

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 74119100
$ws.Range("I122").Value = 90001340
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 270004020
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -270001570
$ws.Range("N122").Value = -10799.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2098.3
$ws.Range("I132").Value = 1175
$ws.Range("J132").Value = 2713.8333
$ws.Range("K132").Value = 3525
$ws.Range("L132").Value = 8141.499899999999
$ws.Range("M132").Value = -995
$ws.Range("N132").Value = -13201.4999
